$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.778.94"
$ws.Range("E2").Value = "  -1.23%  "
$ws.Range("D3").Value = "1.597.44"
$ws.Range("E3").Value = "  -2.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.58"
$ws.Range("E5").Value = "  -2.53%  "
$ws.Range("E6").Value = "  -1.84%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -1.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0617"
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.62"
$ws.Range("E10").Value = "  -2.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0838"
$ws.Range("E11").Value = "  -1.16%  "
$ws.Range("D12").Value = "1.821.05"
$ws.Range("E12").Value = "  -2.05%  "
$ws.Range("D13").Value = "1.613.02"
$ws.Range("E13").Value = "  -1.22%  "
$ws.Range("E14").Value = "  -1.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.529"
$ws.Range("E15").Value = "  -2.20%  "
$ws.Range("D16").Value = "26.741.76"
$ws.Range("E16").Value = "  -1.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.43"
$ws.Range("E17").Value = "  -3.39%  "
$ws.Range("D18").Value = "0.0₃0728"
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "208.91"
$ws.Range("E19").Value = "  -2.33%  "
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("E21").Value = "  -1.44%  "
$ws.Range("E22").Value = "  -2.48%  "
$ws.Range("E23").Value = "  -6.52%  "
$ws.Range("E24").Value = "  -2.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.34"
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.50"
$ws.Range("E26").Value = "  +1.47%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.113"
$ws.Range("E28").Value = "  -4.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.28"
$ws.Range("E29").Value = "  -1.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0500"
$ws.Range("E30").Value = "  -0.66%  "
$ws.Range("E31").Value = "  -2.45%  "
$ws.Range("E32").Value = "  -3.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.668"
$ws.Range("E33").Value = "  +23.68%  "
$ws.Range("E34").Value = "  -1.98%  "
$ws.Range("D35").Value = "1.312.77"
$ws.Range("E35").Value = "  +0.89%  "
$ws.Range("E36").Value = "  -2.69%  "
$ws.Range("E37").Value = "  -0.68%  "
$ws.Range("E38").Value = "  -1.11%  "
$ws.Range("E39").Value = "  -2.66%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("E41").Value = "  -2.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.17"
$ws.Range("E42").Value = "  -3.89%  "
$ws.Range("E43").Value = "  -0.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.84"
$ws.Range("E44").Value = "  +0.90%  "
$ws.Range("D45").Value = "1.734.67"
$ws.Range("E45").Value = "  -1.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.97"
$ws.Range("E46").Value = "  -1.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.61"
$ws.Range("E47").Value = "  +0.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.812"
$ws.Range("E48").Value = "  -0.43%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0102"
$ws.Range("E49").Value = "  -4.57%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0509"
$ws.Range("E50").Value = "  -0.83%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0974"
$ws.Range("E51").Value = "  +2.78%  "
